$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27 (existing rows 27-74 shift down to 28-75).
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
# Most fields mirror the row directly below it (old row 27, now row 28),
# only the date and the price columns (K, L, M, P) differ.
$ws.Range("A27").Value = 9
$ws.Range("B27").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44790
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 100112029
$ws.Range("G27").Value = "Orégano"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 16
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = 20000
$ws.Range("N27").Value = "$/docena de atados"
$ws.Range("O27").Value = "Región Metropolitana"
$ws.Range("P27").Value = 6667
$ws.Range("Q27").Value = 3
$ws.Range("R27").Value = "Hortaliza"
